# Update cryptos list values (Price / Volume(1h)) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is a plain decimal number need to be
# forced to Text format first, otherwise Excel auto-converts them to a
# numeric value (the sheet stores prices as plain text strings).
foreach ($addr in @("D5","D6","D8","D9","D14","D19","D20","D21","D23","D28","D29","D31","D33","D35","D37","D38","D39","D41","D42","D43","D44","D48")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.017.42"
$ws.Range("E2").Value = "  -1.47%  "

$ws.Range("D3").Value = "3.420.08"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "579.43"

$ws.Range("D6").Value = "154.40"
$ws.Range("E6").Value = "  +5.33%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "0.485"
$ws.Range("E8").Value = "  +1.54%  "

$ws.Range("D9").Value = "8.05"
$ws.Range("E9").Value = "  +3.47%  "

$ws.Range("E10").Value = "  +0.84%  "

$ws.Range("E11").Value = "  +3.66%  "

$ws.Range("D12").Value = "4.003.14"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("E13").Value = "  +0.78%  "

$ws.Range("D14").Value = "28.67"
$ws.Range("E14").Value = "  -0.91%  "

$ws.Range("D16").Value = "3.418.75"
$ws.Range("E16").Value = "  -0.49%  "

$ws.Range("D17").Value = "62.031.82"
$ws.Range("E17").Value = "  -1.49%  "

$ws.Range("E18").Value = "  +3.01%  "

$ws.Range("D19").Value = "14.45"
$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("D20").Value = "8.97"
$ws.Range("E20").Value = "  -2.27%  "

$ws.Range("D21").Value = "382.78"
$ws.Range("E21").Value = "  -0.72%  "

$ws.Range("E22").Value = "  +1.62%  "

$ws.Range("D23").Value = "76.04"
$ws.Range("E23").Value = "  +2.13%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").Value = "3.561.09"
$ws.Range("E25").Value = "  -0.82%  "

$ws.Range("E26").Value = "  -1.23%  "

$ws.Range("E27").Value = "  -1.11%  "

$ws.Range("D28").Value = "7.66"
$ws.Range("E28").Value = "  +1.12%  "

$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("E30").Value = "  +1.13%  "

$ws.Range("D31").Value = "7.89"
$ws.Range("E31").Value = "  -2.41%  "

$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("D33").Value = "23.28"
$ws.Range("E33").Value = "  +0.15%  "

$ws.Range("E34").Value = "  +1.38%  "

$ws.Range("D35").Value = "5.58"
$ws.Range("E35").Value = "  +5.75%  "

$ws.Range("E36").Value = "  +0.76%  "

$ws.Range("D37").Value = "6.97"
$ws.Range("E37").Value = "  -1.32%  "

$ws.Range("D38").Value = "168.40"
$ws.Range("E38").Value = "  -0.32%  "

$ws.Range("D39").Value = "31.00"
$ws.Range("E39").Value = "  -1.67%  "

$ws.Range("D40").Value = "3.454.30"
$ws.Range("E40").Value = "  -0.60%  "

$ws.Range("D41").Value = "0.0784"
$ws.Range("E41").Value = "  +2.58%  "

$ws.Range("D42").Value = "42.74"
$ws.Range("E42").Value = "  +0.88%  "

$ws.Range("D43").Value = "0.781"
$ws.Range("E43").Value = "  -0.95%  "

$ws.Range("D44").Value = "4.43"
$ws.Range("E44").Value = "  +1.90%  "

$ws.Range("E45").Value = "  -2.25%  "

$ws.Range("E46").Value = "  -2.32%  "

$ws.Range("D47").Value = "2.551.81"
$ws.Range("E47").Value = "  -0.73%  "

$ws.Range("D48").Value = "23.25"
$ws.Range("E48").Value = "  +2.45%  "

$ws.Range("E49").Value = "  +0.33%  "

$ws.Range("E50").Value = "  -3.17%  "

$ws.Range("E51").Value = "  +0.01%  "
